$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B2").Value = "107号直流"
$ws.Range("C2").Value = 46033.455567129633
$ws.Range("D2").Value = 46034.28329861111
$ws.Range("A3").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B3").Value = "110号直流"
$ws.Range("C3").Value = 46033.563611111109
$ws.Range("D3").Value = 46034.28329861111
$ws.Range("A4").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B4").Value = "208号直流"
$ws.Range("C4").Value = 46033.582951388889
$ws.Range("D4").Value = 46034.28329861111
$ws.Range("A5").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B5").Value = "301号直流"
$ws.Range("C5").Value = 46033.632372685184
$ws.Range("D5").Value = 46034.28329861111
$ws.Range("A6").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B6").Value = "111号直流"
$ws.Range("C6").Value = 46033.686898148146
$ws.Range("D6").Value = 46034.28329861111
$ws.Range("A7").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B7").Value = "206号直流"
$ws.Range("C7").Value = 46033.697546296295
$ws.Range("D7").Value = 46034.28329861111
$ws.Range("A8").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B8").Value = "209号直流"
$ws.Range("C8").Value = 46033.720983796295
$ws.Range("D8").Value = 46034.28329861111
$ws.Range("A9").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B9").Value = "108号直流"
$ws.Range("C9").Value = 46033.733495370368
$ws.Range("D9").Value = 46034.28329861111
$ws.Range("A10").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B10").Value = "101号直流"
$ws.Range("C10").Value = 46033.750289351854
$ws.Range("D10").Value = 46034.28329861111
$ws.Range("A11").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B11").Value = "210号直流"
$ws.Range("C11").Value = 46033.75472222222
$ws.Range("D11").Value = 46034.28329861111
$ws.Range("A12").Value = "飞狐四方坪东区充电站"
$ws.Range("B12").Value = "9176699442100801"
$ws.Range("C12").Value = 46030.706087962964
$ws.Range("D12").Value = 46034.297037037039
$ws.Range("A13").Value = "飞狐四方坪西区充电站"
$ws.Range("B13").Value = "9176699400501304"
$ws.Range("C13").Value = 46032.57640046296
$ws.Range("D13").Value = 46034.297037037039
$ws.Range("A14").Value = "飞狐四方坪西区充电站"
$ws.Range("B14").Value = "9176699400500304"
$ws.Range("C14").Value = 46033.078541666669
$ws.Range("D14").Value = 46034.297037037039
$ws.Range("A15").Value = "飞狐四方坪西区充电站"
$ws.Range("B15").Value = "9176699400501101"
$ws.Range("C15").Value = 46033.224687499998
$ws.Range("D15").Value = 46034.297037037039
$ws.Range("A16").Value = "飞狐四方坪西区充电站"
$ws.Range("B16").Value = "9176699400501303"
$ws.Range("C16").Value = 46033.534456018519
$ws.Range("D16").Value = 46034.297037037039
$ws.Range("A17").Value = "飞狐四方坪西区充电站"
$ws.Range("B17").Value = "9176699400501203"
$ws.Range("C17").Value = 46033.540856481479
$ws.Range("D17").Value = 46034.297037037039
$ws.Range("A18").Value = "飞狐四方坪西区充电站"
$ws.Range("B18").Value = "9176699400500902"
$ws.Range("C18").Value = 46033.54478009259
$ws.Range("D18").Value = 46034.297037037039
$ws.Range("A19").Value = "飞狐四方坪南区充电站"
$ws.Range("B19").Value = "9176699368200101"
$ws.Range("C19").Value = 46033.552187499998
$ws.Range("D19").Value = 46034.297037037039
$ws.Range("A20").Value = "飞狐四方坪西区充电站"
$ws.Range("B20").Value = "9176699400501302"
$ws.Range("C20").Value = 46033.553148148145
$ws.Range("D20").Value = 46034.297037037039
$ws.Range("A21").Value = "飞狐四方坪西区充电站"
$ws.Range("B21").Value = "9176699400500602"
$ws.Range("C21").Value = 46033.557210648149
$ws.Range("D21").Value = 46034.297037037039
$ws.Range("A22").Value = "飞狐四方坪南区充电站"
$ws.Range("B22").Value = "9176699368200103"
$ws.Range("C22").Value = 46033.55914351852
$ws.Range("D22").Value = 46034.297037037039
$ws.Range("A23").Value = "飞狐四方坪南区充电站"
$ws.Range("B23").Value = "9176699368200306"
$ws.Range("C23").Value = 46033.560196759259
$ws.Range("D23").Value = 46034.297037037039
$ws.Range("A24").Value = "飞狐四方坪东区充电站"
$ws.Range("B24").Value = "9176699425700301"
$ws.Range("C24").Value = 46033.561493055553
$ws.Range("D24").Value = 46034.297037037039
$ws.Range("A25").Value = "飞狐四方坪南区充电站"
$ws.Range("B25").Value = "9176699368200406"
$ws.Range("C25").Value = 46033.56354166667
$ws.Range("D25").Value = 46034.297037037039
$ws.Range("A26").Value = "飞狐四方坪东区充电站"
$ws.Range("B26").Value = "9176699442100302"
$ws.Range("C26").Value = 46033.563923611109
$ws.Range("D26").Value = 46034.297037037039
$ws.Range("A27").Value = "飞狐四方坪东区充电站"
$ws.Range("B27").Value = "9176699442100101"
$ws.Range("C27").Value = 46033.568078703705
$ws.Range("D27").Value = 46034.297037037039
$ws.Range("A28").Value = "飞狐四方坪西区充电站"
$ws.Range("B28").Value = "9176699400500303"
$ws.Range("C28").Value = 46033.572245370371
$ws.Range("D28").Value = 46034.297037037039
$ws.Range("A29").Value = "飞狐四方坪西区充电站"
$ws.Range("B29").Value = "9176699400500605"
$ws.Range("C29").Value = 46033.584664351853
$ws.Range("D29").Value = 46034.297037037039
$ws.Range("A30").Value = "飞狐四方坪西区充电站"
$ws.Range("B30").Value = "9176699400500604"
$ws.Range("C30").Value = 46033.586458333331
$ws.Range("D30").Value = 46034.297037037039
$ws.Range("A31").Value = "飞狐四方坪西区充电站"
$ws.Range("B31").Value = "9176699400500102"
$ws.Range("C31").Value = 46033.590543981481
$ws.Range("D31").Value = 46034.297037037039
$ws.Range("A32").Value = "飞狐四方坪西区充电站"
$ws.Range("B32").Value = "9176699400500104"
$ws.Range("C32").Value = 46033.626446759263
$ws.Range("D32").Value = 46034.297037037039
$ws.Range("A33").Value = "飞狐四方坪西区充电站"
$ws.Range("B33").Value = "9176699400500302"
$ws.Range("C33").Value = 46033.635567129626
$ws.Range("D33").Value = 46034.297037037039
$ws.Range("A34").Value = "飞狐四方坪东区充电站"
$ws.Range("B34").Value = "9176699425700302"
$ws.Range("C34").Value = 46033.648958333331
$ws.Range("D34").Value = 46034.297037037039
$ws.Range("A35").Value = "飞狐四方坪西区充电站"
$ws.Range("B35").Value = "9176699400500501"
$ws.Range("C35").Value = 46033.666574074072
$ws.Range("D35").Value = 46034.297037037039
$ws.Range("A36").Value = "飞狐四方坪西区充电站"
$ws.Range("B36").Value = "9176699400501202"
$ws.Range("C36").Value = 46033.708819444444
$ws.Range("D36").Value = 46034.297037037039

$ws.Range("A12:D36").Select()

Write-Host "applied"
